# Update "basic nursing" script: set C12/C13 values on the
# ND-SECOND-YEAR-SECOND-SEMESTER sheet (previously placeholder "-") and make
# that sheet the active / selected tab instead of ND-FIRST-YEAR-SECOND-SEMESTER.

$wb = $excel.ActiveWorkbook

$wsSecondYearSecondSem = $wb.Worksheets.Item("ND-SECOND-YEAR-SECOND-SEMESTER")
$wsFirstYearSecondSem  = $wb.Worksheets.Item("ND-FIRST-YEAR-SECOND-SEMESTER")

# Fill in the previously-missing credit unit (CU) values.
$wsSecondYearSecondSem.Range("C12").Value = 1
$wsSecondYearSecondSem.Range("C13").Value = 4

# Move the selection / active cell on that sheet.
$wsSecondYearSecondSem.Activate()
$wsSecondYearSecondSem.Range("E12").Select()

# The workbook now opens with ND-SECOND-YEAR-SECOND-SEMESTER as the active tab.
$wb.Windows.Item(1).ActiveSheet = $wsSecondYearSecondSem
